$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.107.33"
$ws.Range("E2").Value = "  -2.60%  "

# Row 3
$ws.Range("D3").Value = "2.821.03"
$ws.Range("E3").Value = "  -3.68%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.18"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.02%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.85%  "

# Row 9
$ws.Range("D9").Value = "2.816.67"
$ws.Range("E9").Value = "  -3.87%  "

# Row 10
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.91"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.60%  "

# Row 11
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.102"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -6.53%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.347"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.59%  "

# Row 13
$ws.Range("D13").Value = "3.327.58"
$ws.Range("E13").Value = "  -3.29%  "

# Row 14
$ws.Range("E14").Value = "  +1.48%  "

# Row 15
$ws.Range("D15").Value = "59.227.40"
$ws.Range("E15").Value = "  -2.69%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.47"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -6.61%  "

# Row 17
$ws.Range("D17").Value = "2.842.95"
$ws.Range("E17").Value = "  -2.80%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000134"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.38%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.70"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -6.23%  "

# Row 20
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.96"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.34%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.05"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.59%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.74%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.15%  "

# Row 24
$ws.Range("E24").Value = "  -0.84%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.11"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.42%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.426"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.52%  "

# Row 27
$ws.Range("E27").Value = "  -6.34%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("D29").ClearFormats()

# Row 30
$ws.Range("D30").Value = "0.0₃0792"
$ws.Range("E30").Value = "  -10.03%  "

# Row 31
$ws.Range("E31").Value = "  -0.03%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.59"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.99%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.90"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.81%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.38"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.17%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.13"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.46%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.31"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.54%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.888"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -11.69%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.11"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -8.86%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.70"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.87%  "

# Row 40
$ws.Range("D40").Value = "2.204.62"
$ws.Range("E40").Value = "  -6.69%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.625"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.52%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.51"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.96%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.05%  "

# Row 44
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0554"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.33%  "

# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.34"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -9.26%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.06"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -9.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.35"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.60%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0224"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.56%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0883"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.94%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.53"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -8.90%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.35"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.00%  "
